$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 45154
$ws.Range("J2").Value = 500
$ws.Range("K2").Value = 16500
$ws.Range("L2").Value = 17000
$ws.Range("M2").Value = 16750
$ws.Range("P2").Value = 931

$ws.Range("D3").Value = 45194
$ws.Range("J3").Value = 400
$ws.Range("K3").Value = 16500
$ws.Range("M3").Value = 16750
$ws.Range("P3").Value = 931

$ws.Range("D4").Value = 45005

$ws.Range("D5").Value = 45177
$ws.Range("J5").Value = 540
$ws.Range("K5").Value = 16000
$ws.Range("L5").Value = 17000
$ws.Range("M5").Value = 16500
$ws.Range("P5").Value = 917

$ws.Range("D6").Value = 45159
$ws.Range("J6").Value = 400

$ws.Range("D7").Value = 44557
$ws.Range("J7").Value = 400
$ws.Range("K7").Value = 13000
$ws.Range("L7").Value = 14000
$ws.Range("M7").Value = 13500
$ws.Range("P7").Value = 750

$ws.Range("D8").Value = 44964
$ws.Range("J8").Value = 300
$ws.Range("K8").Value = 20000
$ws.Range("L8").Value = 21000
$ws.Range("M8").Value = 20500
$ws.Range("P8").Value = 1139

$ws.Range("D9").Value = 45166
$ws.Range("J9").Value = 200
$ws.Range("K9").Value = 16000
$ws.Range("M9").Value = 16500
$ws.Range("P9").Value = 917

$ws.Range("D10").Value = 44960
$ws.Range("J10").Value = 400
$ws.Range("K10").Value = 19500
$ws.Range("L10").Value = 20000
$ws.Range("M10").Value = 19750
$ws.Range("P10").Value = 1097

$ws.Range("D12").Value = 45230
$ws.Range("J12").Value = 360
$ws.Range("K12").Value = 16000
$ws.Range("L12").Value = 17000
$ws.Range("M12").Value = 16500
$ws.Range("P12").Value = 917

$ws.Range("D13").Value = 45215
$ws.Range("K13").Value = 16000
$ws.Range("L13").Value = 17000
$ws.Range("M13").Value = 16500
$ws.Range("P13").Value = 917

$ws.Range("D14").Value = 45229
$ws.Range("J14").Value = 460

$ws.Range("D15").Value = 45152
$ws.Range("J15").Value = 500
$ws.Range("K15").Value = 16000
$ws.Range("M15").Value = 16500
$ws.Range("P15").Value = 917

$ws.Range("D16").Value = 44984
$ws.Range("J16").Value = 200

$ws.Range("D17").Value = 44998
$ws.Range("J17").Value = 320
$ws.Range("K17").Value = 17000
$ws.Range("L17").Value = 18000
$ws.Range("M17").Value = 17500
$ws.Range("P17").Value = 972

$ws.Range("D18").Value = 45117
$ws.Range("J18").Value = 300
$ws.Range("K18").Value = 17000
$ws.Range("L18").Value = 18000
$ws.Range("M18").Value = 17500
$ws.Range("P18").Value = 972

$ws.Range("D19").Value = 45068

$ws.Range("D20").Value = 44957
$ws.Range("J20").Value = 400
$ws.Range("K20").Value = 21000
$ws.Range("L20").Value = 22000
$ws.Range("M20").Value = 21500
$ws.Range("P20").Value = 1194

$ws.Range("D22").Value = 45222
$ws.Range("J22").Value = 300

$ws.Range("D23").Value = 45180
$ws.Range("J23").Value = 400
$ws.Range("K23").Value = 16500
$ws.Range("M23").Value = 16750
$ws.Range("P23").Value = 931

$ws.Range("D24").Value = 44547
$ws.Range("J24").Value = 200
$ws.Range("K24").Value = 13000
$ws.Range("L24").Value = 14000
$ws.Range("M24").Value = 13500
$ws.Range("P24").Value = 750

$ws.Range("D25").Value = 44977
$ws.Range("J25").Value = 400
$ws.Range("M25").Value = 16750
$ws.Range("P25").Value = 931
